$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update tire_mu notes (C20): mention "nonlinear" tire model
$ws.Range("C20").Value = "unitless, used if nonlinear tire model not availible. Typically 2/3 of experimental hot mu vale"

# 2. Update tire_mu_correction_factor notes (C21): append "(const radius test)"
$ws.Range("C21").Value = "use with tire data to linearly adjust mu values to adjust belt friction to road friction, based on experimental testing (const radius test)"

# 3. Change B21 cell style from "Good" to "Bad" (compliance highlighting)
$ws.Range("B21").Style = "Bad"

# 4. Add new row 62: toe_deflection_from_rear
$ws.Range("B62").Value = 0.1
$ws.Range("B62").Style = "Neutral"
$ws.Range("C62").Value = "deg per 1kN, per wheel, toe deflection from Fy forces, from experimental testing"
$ws.Range("A62").Value = "toe_deflection_from_rear"

# Update view state to reflect scrolling to the new data / selecting the new row
$ws.Range("A62").Select()

$wb.Save()
